$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item("展览"), $wb.Worksheets.Item("全部类型"))) {
    $ws.Range("F5").Value = 21
    $ws.Range("F6").Value = 114
    $ws.Range("F7").Value = 87
    $ws.Range("F8").Value = 458
    $ws.Range("F11").Value = 569
    $ws.Range("F13").Value = 299
    $ws.Range("F18").Value = 8
    $ws.Range("F22").Value = 928
    $ws.Range("F23").Value = 1396
    $ws.Range("F25").Value = 327
    $ws.Range("F26").Value = 183
    $ws.Range("F27").Value = 74
    $ws.Range("F29").Value = 40
    $ws.Range("F30").Value = 85
    $ws.Range("F32").Value = 248
    $ws.Range("F33").Value = 274
    $ws.Range("F34").Value = 1617
    $ws.Range("F37").Value = 159
    $ws.Range("F38").Value = 581
    $ws.Range("F40").Value = 3639
    $ws.Range("F43").Value = 909
}
